$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 117, shifting existing rows 117-120 down to 118-121.
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new weekly record.
$ws.Range("A117").Value = 6
$ws.Range("B117").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C117").Value = "Metropolitana"
$ws.Range("D117").Value = 44509
$ws.Range("E117").Value = 13
$ws.Range("F117").Value = 100112029
$ws.Range("G117").Value = "Orégano"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 34
$ws.Range("K117").Value = 8000
$ws.Range("L117").Value = 9000
$ws.Range("M117").Value = 8441
$ws.Range("N117").Value = "$/docena de atados"
$ws.Range("O117").Value = "Región Metropolitana"
$ws.Range("P117").Value = 2814
$ws.Range("Q117").Value = 3
$ws.Range("R117").Value = "Hortaliza"
